# Automatic update of files.
# Applies a swap of data between rows 5/6 (partial columns) and rows 11/12
# (full differing columns) in the Artfynd worksheet, matching the target
# OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 <-> Row 6 : swap Id / Ost / Nord / Publik kommentar / Observatorer ---
$ws.Range("A5").Value2  = 130965930
$ws.Range("Q5").Value2  = 496938
$ws.Range("R5").Value2  = 6713359
$ws.Range("AC5").Value2 = "Måttlig förekomst . inventering åt vasa vind"
$ws.Range("AX5").Value2 = "Pia Edfors, Enviro Planning"

$ws.Range("A6").Value2  = 130965940
$ws.Range("Q6").Value2  = 496969
$ws.Range("R6").Value2  = 6713529
$ws.Range("AC6").Value2 = "Betydelsefulla förekomster . inventering åt vasa vind"
$ws.Range("AX6").Value2 = "Pia Edfors, Anders Esplund, Enviro Planning"

# --- Row 11 <-> Row 12 : swap the full species record ---
$ws.Range("A11").Value2  = 130965935
$ws.Range("B11").Value2  = 79243
$ws.Range("D11").Value2  = "NT"
$ws.Range("E11").Value2  = 6425
$ws.Range("F11").Value2  = "Garnlav"
$ws.Range("G11").Value2  = "Alectoria sarmentosa"
$ws.Range("H11").Value2  = "(Ach.) Ach."
$ws.Range("Q11").Value2  = 496969
$ws.Range("R11").Value2  = 6713674
$ws.Range("AC11").Value2 = "Måttlig förekomst . inventering åt vasa vind"
$ws.Range("AX11").Value2 = "Pia Edfors, Enviro Planning"

$ws.Range("A12").Value2  = 130965861
$ws.Range("B12").Value2  = 98930
$ws.Range("D12").Value2  = "LC"
$ws.Range("E12").Value2  = 219790
$ws.Range("F12").Value2  = "Fläcknycklar"
$ws.Range("G12").Value2  = "Dactylorhiza maculata"
$ws.Range("H12").Value2  = "(L.) Soó"
$ws.Range("Q12").Value2  = 497138
$ws.Range("R12").Value2  = 6713448
$ws.Range("AC12").Value2 = "Betydelsefulla förekomster . inventering åt vasa vind"
$ws.Range("AX12").Value2 = "Anders Esplund, Pia Edfors, Enviro Planning"
